# Auto-generated edit script: updates cryptos list Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.623.85"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").Value = "2.488.49"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.26"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.08"
$ws.Range("E6").Value = "  +3.52%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.141"
$ws.Range("E9").Value = "  +3.70%  "

$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").Value = "  +2.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.92"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "2.950.18"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.65"
$ws.Range("E14").Value = "  +1.22%  "

$ws.Range("D15").Value = "67.579.41"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").Value = "2.501.72"
$ws.Range("E17").Value = "  +1.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.53"
$ws.Range("E18").Value = "  +1.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.96"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.65"
$ws.Range("E20").Value = "  +0.33%  "

$ws.Range("E21").Value = "  +1.95%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.59"
$ws.Range("E23").Value = "  +2.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.29"
$ws.Range("E24").Value = "  +1.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.76"
$ws.Range("E25").Value = "  -1.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.12"
$ws.Range("E26").Value = "  -1.45%  "

$ws.Range("D27").Value = "2.620.77"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("D29").Value = "0.0₃0901"
$ws.Range("E29").Value = "  +0.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "506.60"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.80"
$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.122"
$ws.Range("E35").Value = "  +5.08%  "

$ws.Range("E36").Value = "  +2.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.31"
$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.33"
$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  +3.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.328"
$ws.Range("E42").Value = "  +0.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.83"
$ws.Range("E43").Value = "  +0.91%  "

$ws.Range("E44").Value = "  +1.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "144.86"
$ws.Range("E45").Value = "  +1.77%  "

$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.514"
$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("E48").Value = "  +2.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.58"
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.586"
$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("E51").Value = "  +0.59%  "
